$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7169166111259846
$ws.Range("C2").Value = -0.769951142914251
$ws.Range("D2").Value = 0.81811524812827

$ws.Range("B3").Value = 0.8247668291515247
$ws.Range("C3").Value = -0.8638716263040771
$ws.Range("D3").Value = 0.7044915131233094

$ws.Range("B4").Value = -0.7316417397939979
$ws.Range("C4").Value = -0.8040945298166184
$ws.Range("D4").Value = 0.7861097633820693

$ws.Range("B5").Value = 0.7557354537191054
$ws.Range("C5").Value = 0.7016553533637517
$ws.Range("D5").Value = 0.5881177338195093

$ws.Range("B6").Value = -0.7833197689211761
$ws.Range("C6").Value = -0.7456114469046508
$ws.Range("D6").Value = 0.5972052172252177

$ws.Range("B7").Value = 0.6812417260515466
$ws.Range("C7").Value = -0.6780949752640816
$ws.Range("D7").Value = -0.5882442160195755

$ws.Range("B8").Value = -0.7604195397640003
$ws.Range("C8").Value = 0.7695151871108438
$ws.Range("D8").Value = -0.6103447749713604

$ws.Range("B9").Value = -0.8024343564222106
$ws.Range("C9").Value = 0.8619570517734312
$ws.Range("D9").Value = 0.8295627460520123
